$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 11 by copying row 10 (same columns A:AQ) so it inherits the
# same cell layout/formatting (incl. the blank column-E placeholder cell)
# as every other data row, then overwrite with the 2021 values.
$ws.Range("A10:AQ10").Copy($ws.Range("A11"))

$ws.Cells.Item(11, 1).Value = "2021年"

$values = @{
    2  = 877.3099999999999
    3  = 250.61
    4  = 26.03
    6  = 415.81
    7  = 1101.02
    8  = 63.17
    9  = 3833.49
    10 = 85.89
    11 = 20914.43
    12 = 180.1
    13 = 9.59
    14 = 1.76
    15 = 182.86
    16 = 228.53
    17 = 11.63
    18 = 56.62
    19 = 385.86
    20 = 129.77
    21 = 1631.93
    22 = 154.36
    23 = 290.87
    24 = 155.15
    25 = 36.27
    26 = 1987.27
    27 = 166.96
    28 = 258.24
    29 = 33.98
    30 = 230.99
    31 = 368.83
    32 = 2520.51
    33 = 886.72
    34 = 265.69
    35 = 1000.66
    36 = 5.37
    37 = 373.94
    38 = 207.46
    39 = 585.87
    40 = 25.07
    41 = 1387.86
    42 = 470.97
    43 = 29
}

foreach ($col in $values.Keys) {
    $ws.Cells.Item(11, $col).Value = $values[$col]
}

$wb.Save()
